# page_22: split the single raw grid into two tables.
#  - Sheet1 becomes "page_22 - Table 1": a 2-line header block (merged
#    A1:Z1 / A2:Z2) carrying the extracted document title + review-sheet
#    caption, followed by the original bordered grid (now at rows 4-7,
#    row 3 left blank as a spacer row between header and grid).
#  - A new "page_22 - Table 2" sheet holds the remaining bordered grid
#    (2 rows) using the same cell formatting as before.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Second sheet, placed after the first, carrying the tail of the old grid.
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws1.Name = "page_22 - Table 1"
$ws2.Name = "page_22 - Table 2"

# --- Table 2: plain bordered 3x2 grid (was rows 6:7 of the old sheet) ---
$t2 = $ws2.Range("A1:C2")
$t2.Borders.LineStyle = 1
$t2.HorizontalAlignment = -4131
$t2.VerticalAlignment = -4160
$t2.WrapText = $true

# --- Table 1: new 2-line text header + the remaining bordered grid rows ---

# Row 3 (former row 3 of the old grid) becomes an empty spacer row.
$ws1.Range("A3:C3").Clear()

# Rows 4 and 5 keep the bordered-grid look (row 5 is the old row 4 shifted
# down to make room for row 3 as a spacer).
$gridRows = $ws1.Range("A4:C5")
$gridRows.Borders.LineStyle = 1
$gridRows.HorizontalAlignment = -4131
$gridRows.VerticalAlignment = -4160
$gridRows.WrapText = $true

# Header line 1: " MASTER PACKAGE"
$h1 = $ws1.Range("A1:Z1")
$h1.ClearFormats()
$h1.Font.Size = 13
$h1.HorizontalAlignment = -4131
$h1.WrapText = $true
$ws1.Range("A1").Value = " MASTER PACKAGE"
$h1.Merge()

# Header line 2: long bold review-sheet caption
$h2 = $ws1.Range("A2:Z2")
$h2.Font.Size = 13
$h2.Font.Bold = $true
$h2.HorizontalAlignment = -4131
$h2.WrapText = $true
$ws1.Range("A2").Value = " WesternGlove Centric8 PROD                                         M12225BVS563:KONRAD                                     DUP REVIEW SHEET                                 MASTER"
$h2.Merge()

Write-Host "done"
